$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 677
$ws.Range("F9").Value = 7343
$ws.Range("F11").Value = 148
$ws.Range("F13").Value = 34
$ws.Range("F17").Value = 1787
$ws.Range("F18").Value = 1082
$ws.Range("F19").Value = 26
$ws.Range("F23").Value = 1244
$ws.Range("F26").Value = 1127
$ws.Range("F29").Value = 130
$ws.Range("F31").Value = 3191
$ws.Range("F32").Value = 2266
$ws.Range("F33").Value = 3888
$ws.Range("F37").Value = 1152
$ws.Range("F42").Value = 156
$ws.Range("F43").Value = 521

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 94
$ws.Range("F20").Value = 60

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F12").Value = 424
$ws.Range("F13").Value = 1827
$ws.Range("F14").Value = 8176

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 677
$ws.Range("F11").Value = 424
$ws.Range("F13").Value = 34
$ws.Range("F17").Value = 1082
$ws.Range("F18").Value = 26
$ws.Range("F22").Value = 1244
$ws.Range("F25").Value = 1127
$ws.Range("F26").Value = 94
$ws.Range("F31").Value = 60
$ws.Range("F32").Value = 130
$ws.Range("F34").Value = 3191
$ws.Range("F35").Value = 2266
$ws.Range("F36").Value = 3888
$ws.Range("F40").Value = 1152
$ws.Range("F43").Value = 156
$ws.Range("F45").Value = 521
